$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "_old" / "_new" header label suffixes to "_FV2404" / "_FV2410"
#    (these occur only in the row-1 header cells, 10 occurrences each).
$ws.Cells.Replace("_old", "_FV2404") | Out-Null
$ws.Cells.Replace("_new", "_FV2410") | Out-Null

# 2. Turn the full used range into an Excel Table ("Table1") with an AutoFilter,
#    matching the workbook's data extent (A1:U79).
$rng = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# 3. Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
